# Rename the inline picture shapes living in the document's headers/footers.
#
# The Pearson logo pictures (shown in the footers) are relabeled from
# "image2.png" to "image1.png", and the BTec logo picture (shown in the
# header) is relabeled from "image1.jpg" to "image2.jpg". The pictures'
# actual embedded data/relationships are untouched - only their display
# name (wp:docPr/@name, surfaced here as InlineShape.Name) changes.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    for ($hi = 1; $hi -le 3; $hi++) {
        $hdr = $sec.Headers.Item($hi)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($ii = 1; $ii -le $shapes.Count; $ii++) {
                $shp = $shapes.Item($ii)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    for ($fi = 1; $fi -le 3; $fi++) {
        $ftr = $sec.Footers.Item($fi)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($ii = 1; $ii -le $shapes.Count; $ii++) {
                $shp = $shapes.Item($ii)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
